$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.035578370094299
$ws.Range("B1").Value = 1.248359799385071
$ws.Range("C1").Value = 1.674901604652405
$ws.Range("D1").Value = 3.227506160736084
$ws.Range("E1").Value = 2.499902248382568
